$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3719
$ws.Range("B2").Value = 6321
$ws.Range("C2").Value = 2347
$ws.Range("D2").Value = 6545
$ws.Range("E2").Value = 2064
$ws.Range("F2").Value = 7226
